$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: Закуски / Рыба / 10000
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("D10").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Rows.Item(11).RowHeight = $ws.Rows.Item(10).RowHeight
$ws.Range("A11").Value = "Закуски"
$ws.Range("B11").Value = "Рыба"
$ws.Range("D11").Value = 10000

# Row 12: Закуски / Мясо / 20000
$ws.Range("A10").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("D10").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Rows.Item(12).RowHeight = $ws.Rows.Item(10).RowHeight
$ws.Range("A12").Value = "Закуски"
$ws.Range("B12").Value = "Мясо"
$ws.Range("D12").Value = 20000

# Row 13: Закуски / Сыр / 40000
$ws.Range("A10").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("D10").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Rows.Item(13).RowHeight = $ws.Rows.Item(10).RowHeight
$ws.Range("A13").Value = "Закуски"
$ws.Range("B13").Value = "Сыр"
$ws.Range("D13").Value = 40000

$ws.Range("D14").Select()
